# Fix the typo'd worksheet name and make "Course Equivalents" the active sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Course Equivelents")
$ws.Name = "Course Equivalents"

$ws.Activate()
